$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14-18 down to 15-19
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new data record
$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 45089
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100107
$ws.Cells.Item(14, 8).Value = "Otros"
$ws.Cells.Item(14, 9).Value = 100107011
$ws.Cells.Item(14, 10).Value = "Tuna"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 14).Value = 22000
$ws.Cells.Item(14, 15).Value = 23000
$ws.Cells.Item(14, 16).Value = 22500
$ws.Cells.Item(14, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(14, 18).Value = "Región Metropolitana"
$ws.Cells.Item(14, 19).Value = 1250
$ws.Cells.Item(14, 20).Value = 18
